$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.469.08'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '2.270.25'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '120.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +7.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.59'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('E7').Value = '  +3.08%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.54'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0942'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('E12').Value = '  +5.06%  '
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.77'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.57%  '
$ws.Range('E15').Value = '  +5.90%  '
$ws.Range('D16').Value = '2.610.55'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').Value = '2.266.63'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '43.621.08'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.90'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.66'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('E22').Value = '  -4.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.98'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.74'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.96'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.32'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +8.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.01'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.33%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.97'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0918'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.74'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('E35').Value = '  +14.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.130'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0383'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +7.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.73'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.04%  '
$ws.Range('E39').Value = '  +4.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.57'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.244'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.81'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.62'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -8.21%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.73'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -8.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '77.27'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +37.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.667'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +18.78%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.27'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.59'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.26%  '
$ws.Range('E51').Value = '  +1.61%  '
